# Update mapa_interactivo.xlsx:
#  - Fill in previously-missing geocoded coordinates (M/N) and
#    Operacion/Zona (O/P) for two already-existing rows (cases 6357 and
#    -502) which appear in both "General" and "Optical_Power".
#  - Append 6 new case rows (6362, 6363, 6372, 6376, 6377, 6383) to the
#    "General" sheet (which holds every case), and also append each one
#    to its provider-specific sheet (NEW, AYKO, PEBCOM, Optical_Power).
#
# NOTE: columns A-L, O, P are plain text in this workbook (even when the
# text looks numeric, e.g. case ids or "12"), only columns M/N
# (Coordenada_X / Coordenada_Y) are real numbers. Excel's COM layer will
# silently reinterpret numeric-looking or date-looking strings (e.g.
# "6362", "7/8/2025") unless the cell is explicitly forced to text
# format first; the cell is then reset back to the default/general
# format so no stray styling is left behind.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

function Set-NumCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Writes a full A..P case row (12 text cols, 2 numeric coord cols, 2 text cols).
function Set-CaseRow($ws, $row, $vals) {
    Set-TextCell $ws $row 1  $vals[0]
    Set-TextCell $ws $row 2  $vals[1]
    Set-TextCell $ws $row 3  $vals[2]
    Set-TextCell $ws $row 4  $vals[3]
    Set-TextCell $ws $row 5  $vals[4]
    Set-TextCell $ws $row 6  $vals[5]
    Set-TextCell $ws $row 7  $vals[6]
    Set-TextCell $ws $row 8  $vals[7]
    Set-TextCell $ws $row 9  $vals[8]
    Set-TextCell $ws $row 10 $vals[9]
    Set-TextCell $ws $row 11 $vals[10]
    Set-TextCell $ws $row 12 $vals[11]
    Set-NumCell  $ws $row 13 $vals[12]
    Set-NumCell  $ws $row 14 $vals[13]
    Set-TextCell $ws $row 15 $vals[14]
    Set-TextCell $ws $row 16 $vals[15]
}

# ---------------------------------------------------------------------
# New case data (columns A..P)
# ---------------------------------------------------------------------
$case6362 = @("6362", "7/8/2025", "ARIAS 3422", "12", "808099435", "NEW", `
    "Pendiente", "Poste inclinado mal ubicado", "1", "Cambio", "Sin equipos", `
    "Poste", -58.483313, -34.54605, "Saavedra", "Capital Norte")

$case6363 = @("6363", "7/8/2025", "MOLDES 3730", "12", "808099415", "NEW", `
    "Pendiente", "Poste inclinado", "1", "Aplomo", "Sin equipos", `
    "Poste", -58.47192, -34.549398, "Saavedra", "Capital Norte")

$case6372 = @("6372", "7/8/2025", "AVELLANEDA 4500", "10", "808099405", "AYKO", `
    "Pendiente", "Picada", "1", "Cambio", "Sin equipos", `
    "Pasante", -58.489219, -34.632475, "Devoto", "Capital Norte")

$case6376 = @("6376", "7/8/2025", "BOYACA 712", "7", "808099366", "PEBCOM", `
    "Pendiente", "Picada", "1", "Cambio", "Sin equipos", `
    "Pasante", -58.461858, -34.619348, "Boedo", "Capital Sur")

$case6377 = @("6377", "7/8/2025", "GUARDIA VIEJA 4377", "5", "808099347", "Optical Power", `
    "Pendiente", "Picada", "1", "Cambio", "Sin equipos", `
    "Pasante", -58.426322, -34.600097, "Almagro", "Capital Sur")

$case6383 = @("6383", "7/8/2025", "FALCON, RAMON L.,CNEL. 1411", "6", "808099320", "Optical Power", `
    "Pendiente", "Picada", "1", "Cambio", "Sin equipos", `
    "Pasante", -58.448523, -34.62452, "Boedo", "Capital Sur")

# ---------------------------------------------------------------------
# General sheet: fill in missing coordinates for rows 387/388, then
# append the 6 new cases as rows 389-394.
# ---------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("General")

Set-NumCell  $wsGeneral 387 13 -58.473179
Set-NumCell  $wsGeneral 387 14 -34.629138
Set-TextCell $wsGeneral 387 15 "Devoto"
Set-TextCell $wsGeneral 387 16 "Capital Norte"

Set-NumCell  $wsGeneral 388 13 -58.400188
Set-NumCell  $wsGeneral 388 14 -34.583882
Set-TextCell $wsGeneral 388 15 "Recoleta"
Set-TextCell $wsGeneral 388 16 "Capital Sur"

Set-CaseRow $wsGeneral 389 $case6362
Set-CaseRow $wsGeneral 390 $case6363
Set-CaseRow $wsGeneral 391 $case6372
Set-CaseRow $wsGeneral 392 $case6376
Set-CaseRow $wsGeneral 393 $case6377
Set-CaseRow $wsGeneral 394 $case6383

# ---------------------------------------------------------------------
# PEBCOM sheet: append case 6376 as row 49.
# ---------------------------------------------------------------------
$wsPebcom = $wb.Worksheets.Item("PEBCOM")
Set-CaseRow $wsPebcom 49 $case6376

# ---------------------------------------------------------------------
# Optical_Power sheet: fill in missing coordinates for rows 76/77, then
# append cases 6377 and 6383 as rows 78-79.
# ---------------------------------------------------------------------
$wsOptical = $wb.Worksheets.Item("Optical_Power")

Set-NumCell  $wsOptical 76 13 -58.473179
Set-NumCell  $wsOptical 76 14 -34.629138
Set-TextCell $wsOptical 76 15 "Devoto"
Set-TextCell $wsOptical 76 16 "Capital Norte"

Set-NumCell  $wsOptical 77 13 -58.400188
Set-NumCell  $wsOptical 77 14 -34.583882
Set-TextCell $wsOptical 77 15 "Recoleta"
Set-TextCell $wsOptical 77 16 "Capital Sur"

Set-CaseRow $wsOptical 78 $case6377
Set-CaseRow $wsOptical 79 $case6383

# ---------------------------------------------------------------------
# NEW sheet: append cases 6362 and 6363 as rows 53-54.
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Item("NEW")
Set-CaseRow $wsNew 53 $case6362
Set-CaseRow $wsNew 54 $case6363

# ---------------------------------------------------------------------
# AYKO sheet: append case 6372 as row 112.
# ---------------------------------------------------------------------
$wsAyko = $wb.Worksheets.Item("AYKO")
Set-CaseRow $wsAyko 112 $case6372
